# Daylight Savings Time update: rename event timer rows from explicit
# clock times (e.g. "12PM", "830PM", "10PM") to generic sequential
# occurrence labels (e.g. "1", "2", "3"), and re-order / expand the
# Wrathborne Invasion entries into two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrathborne Invasion now has two occurrences (was a single 12PM row).
$ws.Range("A12").Value = "Wrathborne Invasion 1"
$ws.Range("A13").Value = "Wrathborne Invasion 2"

# Ancient Nightmare occurrences.
$ws.Range("A2").Value = "Ancient Nightmare 1"
$ws.Range("A3").Value = "Ancient Nightmare 2"
$ws.Range("A4").Value = "Ancient Nightmare 3"

# Haunted Carriage occurrences.
$ws.Range("A5").Value = "Haunted Carriage 1"
$ws.Range("A6").Value = "Haunted Carriage 2"
$ws.Range("A7").Value = "Haunted Carriage 3"

# Demon Gates occurrences.
$ws.Range("A8").Value = "Demon Gates 1"
$ws.Range("A9").Value = "Demon Gates 2"
$ws.Range("A10").Value = "Demon Gates 3"

# Ancient Arena no longer has a fixed time suffix.
$ws.Range("A11").Value = "Ancient Arena"

# Move the active selection to match the saved view state.
$ws.Range("F9").Select()
